# Generate Report for Handoff
# Adds a new tracked file (bb0cc955-f27b-4958-90e3-615021de74f4.md) as row 3
# to the Overview / zh-cn / de-de worksheets, mirroring the existing row 2
# (962208f6-d743-44a1-b690-3777c6ada601.md) pattern.
#
# Note: a leading "'" forces a literal-text cell (matches how the source
# workbook stores "True"/"False"/"" as shared strings, not booleans/blanks).

$wb = $excel.ActiveWorkbook

$newFileBase  = "bb0cc955-f27b-4958-90e3-615021de74f4"
$newMdName    = "$newFileBase.md"
$newMdPath    = "e2e\$newFileBase.md"
$newMdUrl     = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c5b3ed7bce5f8ae1c45dc2c0a14f9baf834ccee0/e2e/$newMdName"

$zhXlfName    = "$newFileBase.a0aa9eadc0f31ad8225b74c28297a67796c77ea5.zh-cn.xlf"
$deXlfName    = "$newFileBase.a0aa9eadc0f31ad8225b74c28297a67796c77ea5.de-de.xlf"

$overviewDate  = "2016-08-25 06:38:44"
$zhHandoffDate = "2016-08-25 06:38:39"
$deHandoffDate = "2016-08-25 06:38:44"
$epochDate     = "0001-01-01 00:00:00"

$dateFormat = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------------
# Sheet "Overview" -> new row 3
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A3").Value = $newMdName
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("D3").Value = "'"
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = $overviewDate
$wsOverview.Range("G3").NumberFormat = $dateFormat

$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $newMdUrl, $null, $null, $newMdPath)

$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.Resize($wsOverview.Range("A1:G3"))

# ---------------------------------------------------------------------------
# Sheet "zh-cn" -> new row 3
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("B3").Value = ".md"
$wsZh.Range("C3").Value = "Ready for handoff"
$wsZh.Range("D3").Value = "e2e"
$wsZh.Range("E3").Value = "ht"
$wsZh.Range("F3").Value = "'False"
$wsZh.Range("G3").Value = $zhXlfName
$wsZh.Range("H3").Value = $zhHandoffDate
$wsZh.Range("H3").NumberFormat = $dateFormat
$wsZh.Range("I3").Value = "'"
$wsZh.Range("J3").Value = "'"
$wsZh.Range("K3").Value = $epochDate
$wsZh.Range("K3").NumberFormat = $dateFormat
$wsZh.Range("L3").Value = "'"
$wsZh.Range("M3").Value = "'True"
$wsZh.Range("N3").Value = "'"
$wsZh.Range("O3").Value = "'False"
$wsZh.Range("P3").Value = "'"

$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $newMdUrl, $null, $null, $newMdName)

$loZh = $wsZh.ListObjects.Item(1)
$loZh.Resize($wsZh.Range("A1:P3"))

# ---------------------------------------------------------------------------
# Sheet "de-de" -> new row 3
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("B3").Value = ".md"
$wsDe.Range("C3").Value = "Ready for handoff"
$wsDe.Range("D3").Value = "e2e"
$wsDe.Range("E3").Value = "ht"
$wsDe.Range("F3").Value = "'False"
$wsDe.Range("G3").Value = $deXlfName
$wsDe.Range("H3").Value = $deHandoffDate
$wsDe.Range("H3").NumberFormat = $dateFormat
$wsDe.Range("I3").Value = "'"
$wsDe.Range("J3").Value = "'"
$wsDe.Range("K3").Value = $epochDate
$wsDe.Range("K3").NumberFormat = $dateFormat
$wsDe.Range("L3").Value = "'"
$wsDe.Range("M3").Value = "'True"
$wsDe.Range("N3").Value = "'"
$wsDe.Range("O3").Value = "'False"
$wsDe.Range("P3").Value = "'"

$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $newMdUrl, $null, $null, $newMdName)

$loDe = $wsDe.ListObjects.Item(1)
$loDe.Resize($wsDe.Range("A1:P3"))
